# Adds 4 new game result rows (60-63) to Sheet1, matching the new
# "Data/bombay1.xlsx" play sessions recorded after the last existing row (59).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 60; A = "2024-06-06 18:52:18"; B = 39; C = 28; D = 3;  E = 10; F = 4; G = 11; H = 0; I = 0.001; J = 0.05; K = 0.003; L = 100; M = 500; N = 10; O = 9; P = 3; Q = 1000; R = 5; S = 2; T = 70; U = 0.717948717948718;  V = "Data/bombay1.xlsx"; W = -103000; X = "No es Simulación" },
    @{ Row = 61; A = "2024-06-06 23:24:27"; B = 39; C = 29; D = 2;  E = 15; F = 6; G = 6;  H = 0; I = 0.001; J = 0.05; K = 0.003; L = 100; M = 500; N = 10; O = 9; P = 3; Q = 1000; R = 5; S = 2; T = 70; U = 0.7435897435897436; V = "Data/bombay1.xlsx"; W = -29000;  X = "No es Simulación" },
    @{ Row = 62; A = "2024-06-10 22:38:38"; B = 14; C = 10; D = 4;  E = 2;  F = 1; G = 3;  H = 0; I = 0.001; J = 0.05; K = 0.003; L = 100; M = 500; N = 10; O = 9; P = 3; Q = 1000; R = 5; S = 2; T = 70; U = 0.7142857142857143; V = "Data/bombay1.xlsx"; W = 150000;  X = "No es Simulación" },
    @{ Row = 63; A = "2024-06-13 19:36:58"; B = 35; C = 21; D = 2;  E = 5;  F = 5; G = 9;  H = 0; I = 0.001; J = 0.05; K = 0.003; L = 100; M = 500; N = 10; O = 9; P = 3; Q = 1000; R = 5; S = 2; T = 70; U = 0.6;                 V = "Data/bombay1.xlsx"; W = -907000; X = "No es Simulación" }
)

foreach ($row in $newRows) {
    $n = $row.Row
    $ws.Range("A$n").Value = $row.A
    $ws.Range("B$n").Value = $row.B
    $ws.Range("C$n").Value = $row.C
    $ws.Range("D$n").Value = $row.D
    $ws.Range("E$n").Value = $row.E
    $ws.Range("F$n").Value = $row.F
    $ws.Range("G$n").Value = $row.G
    $ws.Range("H$n").Value = $row.H
    $ws.Range("I$n").Value = $row.I
    $ws.Range("J$n").Value = $row.J
    $ws.Range("K$n").Value = $row.K
    $ws.Range("L$n").Value = $row.L
    $ws.Range("M$n").Value = $row.M
    $ws.Range("N$n").Value = $row.N
    $ws.Range("O$n").Value = $row.O
    $ws.Range("P$n").Value = $row.P
    $ws.Range("Q$n").Value = $row.Q
    $ws.Range("R$n").Value = $row.R
    $ws.Range("S$n").Value = $row.S
    $ws.Range("T$n").Value = $row.T
    $ws.Range("U$n").Value = $row.U
    $ws.Range("V$n").Value = $row.V
    $ws.Range("W$n").Value = $row.W
    $ws.Range("X$n").Value = $row.X
}
